$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B13").Value = "colabfold"
$ws.Range("A13").Value = "Notebook"
$ws.Range("C13").Value = "1.5.5"
$ws.Range("D13").Value = "GUI to alphafold"
